# "D9 ->" sheet is the active/tab-selected sheet in this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordnance inventory update for D1.2: record ordnance expended this period
# in column D for the rows that changed.
$ws.Range("D7").Value = 3
$ws.Range("D10").Value = 1
$ws.Range("D13").Value = 4
$ws.Range("D18").Value = 14
$ws.Range("D19").Value = 6

# Recompute the "Current inventory" column (T) for every ordnance row as one
# shared formula (mirrors Excel's own fill-down behaviour, which collapses
# the identical relative formula typed/filled across T4:T19 into a single
# shared formula group).
$ws.Range("T4:T19").Formula = "=SUM(C4:C4)-SUM(D4:S4)"

# Header label moves from "D0" to "D1.1" (current-as-of marker).
$ws.Range("T1").Value = "D1.1"

# Leave the cursor where the author left it when saving.
$ws.Range("I12").Select() | Out-Null
